$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "64.938.45"
$ws.Range("D3").Value = "3.139.84"
$ws.Range("E3").Value = "  +2.09%  "
$ws.Range("E4").Value = "  -0.01%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "570.45"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  +2.01%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "150.56"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  +4.10%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "3.136.24"
$ws.Range("E8").Value = "  +2.15%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.525"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  +4.11%  "
$ws.Range("E10").Value = "  +6.47%  "
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "6.13"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -1.48%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.502"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  +6.95%  "
$ws.Range("E13").Value = "  +12.25%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "37.42"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  +6.69%  "
$ws.Range("D15").Value = "3.659.35"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "64.980.05"
$ws.Range("E16").Value = "  +1.11%  "
$ws.Range("E17").Value = "  +6.33%  "
$ws.Range("D18").Value = "3.151.79"
$ws.Range("E18").Value = "  +2.47%  "
$ws.Range("E19").Value = "  +0.45%  "
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = "510.42"
$cell.Style = "Normal"
$ws.Range("E20").Value = "  +6.95%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "14.86"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  +6.91%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "0.728"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  +8.18%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "15.53"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +13.72%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "7.82"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  +3.66%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "85.43"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +5.21%  "
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +4.47%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "8.71"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  +8.63%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "2.18"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  +5.48%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "27.90"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  +6.91%  "
$ws.Range("E31").Value = "  -0.08%  "
$ws.Range("E32").Value = "  +4.60%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "2.64"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  +6.13%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "6.01"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  +8.17%  "
$ws.Range("E35").Value = "  +6.38%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "55.54"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.47%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "476.01"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +4.27%  "
$ws.Range("E38").Value = "  +4.37%  "
$ws.Range("E39").Value = "  +4.17%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "3.02"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  +0.26%  "
$ws.Range("D41").Value = "3.107.50"
$ws.Range("E41").Value = "  +5.17%  "
$ws.Range("E42").Value = "  +4.53%  "
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.119"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  +3.78%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "0.290"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  +11.44%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "2.41"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  +12.68%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "29.05"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  +4.43%  "
$ws.Range("D47").Value = "0.0₃0576"
$ws.Range("E47").Value = "  +11.71%  "
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E49").Value = "  +3.89%  "
$ws.Range("E50").Value = "  +10.54%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "118.85"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -1.78%  "
